$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row (B, C, D, E, G). F column is unchanged by this edit.
$data = @{
    2 = @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.426980108624251)
    3 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5 = @(0.04172184405617529, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.103368080369475)
    6 = @(0.2881169905109251, 1.626987699542094, 3.223369029078222, 2797.565817734744, 2802.704291453875)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
